$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells are treated as text so numeric-looking
# values (e.g. "1.001") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Coin names (column B)
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("B46").Value = 'Aptos'
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("B51").Value = 'Decentraland'

# Links (column C)
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

# Prices (column D)
$ws.Range("D2").Value = '30.590.31'
$ws.Range("D3").Value = '1.928.04'
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").Value = '245.80'
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").Value = '0.4744'
$ws.Range("D8").Value = '0.2904'
$ws.Range("D9").Value = '0.06806'
$ws.Range("D10").Value = '105.97'
$ws.Range("D11").Value = '18.46'
$ws.Range("D12").Value = '0.07767'
$ws.Range("D13").Value = '1.907.04'
$ws.Range("D14").Value = '5.337'
$ws.Range("D15").Value = '0.6703'
$ws.Range("D16").Value = '288.04'
$ws.Range("D17").Value = '30.621.11'
$ws.Range("D18").Value = '0.000007638'
$ws.Range("D19").Value = '13.03'
$ws.Range("D20").Value = '1.001'
$ws.Range("D21").Value = '2.165.32'
$ws.Range("D22").Value = '5.404'
$ws.Range("D23").Value = '1.001'
$ws.Range("D24").Value = '6.272'
$ws.Range("D25").Value = '9.366'
$ws.Range("D26").Value = '168.51'
$ws.Range("D27").Value = '20.94'
$ws.Range("D29").Value = '0.1083'
$ws.Range("D30").Value = '1.370'
$ws.Range("D31").Value = '4.157'
$ws.Range("D32").Value = '4.001'
$ws.Range("D33").Value = '0.05069'
$ws.Range("D34").Value = '0.7375'
$ws.Range("D35").Value = '1.155'
$ws.Range("D36").Value = '0.02092'
$ws.Range("D37").Value = '2.723'
$ws.Range("D38").Value = '2.692'
$ws.Range("D39").Value = '2.065'
$ws.Range("D40").Value = '110.85'
$ws.Range("D41").Value = '0.8763'
$ws.Range("D42").Value = '0.4378'
$ws.Range("D43").Value = '5.906'
$ws.Range("D44").Value = '1.001'
$ws.Range("D45").Value = '67.83'
$ws.Range("D46").Value = '7.252'
$ws.Range("D47").Value = '9.371'
$ws.Range("D48").Value = '48.55'
$ws.Range("D49").Value = '0.1232'
$ws.Range("D50").Value = '35.30'
$ws.Range("D51").Value = '0.4089'

# Volume(1h) percentages (column E)
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("E3").Value = '  +1.21%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("E8").Value = '  -2.28%  '
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("E10").Value = '  +3.80%  '
$ws.Range("E11").Value = '  -4.11%  '
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("E14").Value = '  +2.82%  '
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("E16").Value = '  -7.29%  '
$ws.Range("E17").Value = '  -1.13%  '
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("E19").Value = '  -1.93%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("E22").Value = '  +3.01%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("E25").Value = '  -0.50%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("E28").Value = '  +5.31%  '
$ws.Range("E29").Value = '  -3.33%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("E32").Value = '  -0.59%  '
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("E35").Value = '  -1.73%  '
$ws.Range("E36").Value = '  +5.10%  '
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("E42").Value = '  +2.86%  '
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("E48").Value = '  +13.69%  '
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("E51").Value = '  +5.03%  '

# Restore default (Normal) style on the Price column so no stray number format
# is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
